# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status cells move from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Per-language sheets (zh-cn, de-de) get their "Latest Target File" /
#    "Latest Handback File" (and, where relevant, "Latest Handback DateTime")
#    columns populated, with the target-file cell becoming a hyperlink to the
#    source markdown file (mirroring the existing column A hyperlinks).
#  - Columns that now hold longer text are widened to fit.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb13823943dc1a7e6185be51461bc8a841a8fba1/e2e/"
$file1Name = "4cca8c94-4788-41c4-9ed3-916e45552559.md"
$file2Name = "f6e5b3dd-0781-449c-bc4b-d0cf6f5878fa.md"
$file1Url  = $mdUrlBase + $file1Name
$file2Url  = $mdUrlBase + $file2Name

# Quantized column widths: the stored OOXML column width is
# ColumnWidth + 5/6, rounded to the nearest 1/6 by this host - so pick the
# ColumnWidth that lands closest to the desired stored width.
$wideStatusColWidth = 29.166666666666668   # -> stored ~30   (was ~17.22)
$wideFileColWidth   = 39.166666666666664   # -> stored 40    (was ~18.65 / ~21.71)

# Match the look of the workbook's existing custom "HyperLink" cell style
# (underline + #6495ED) used by column A, rather than the theme-blue style
# that Hyperlinks.Add applies by default.
$hyperlinkColor = 15570276   # OLE BGR for RGB(100,149,237) / #6495ED

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) show the new status
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusColWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $file1Url, $null, $null, $file1Name) | Out-Null
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Range("J2").Value = "4cca8c94-4788-41c4-9ed3-916e45552559.8e5fb013a8ed569c74e17848fbe74211d99e9ae1.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $file2Url, $null, $null, $file2Name) | Out-Null
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = $hyperlinkColor
$wsZh.Range("J3").Value = "f6e5b3dd-0781-449c-bc4b-d0cf6f5878fa.f5e7e4bca0ea7a95afdd34a4cdfd682c2d4ad5ab.zh-cn.xlf"

$wsZh.Columns.Item(3).ColumnWidth = $wideStatusColWidth
$wsZh.Columns.Item(9).ColumnWidth = $wideFileColWidth
$wsZh.Columns.Item(10).ColumnWidth = $wideFileColWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $file1Url, $null, $null, $file1Name) | Out-Null
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Range("J2").Value = "4cca8c94-4788-41c4-9ed3-916e45552559.8e5fb013a8ed569c74e17848fbe74211d99e9ae1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-18 08:30:44"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $file2Url, $null, $null, $file2Name) | Out-Null
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = $hyperlinkColor
$wsDe.Range("J3").Value = "f6e5b3dd-0781-449c-bc4b-d0cf6f5878fa.f5e7e4bca0ea7a95afdd34a4cdfd682c2d4ad5ab.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-18 08:30:44"

$wsDe.Columns.Item(3).ColumnWidth = $wideStatusColWidth
$wsDe.Columns.Item(9).ColumnWidth = $wideFileColWidth
$wsDe.Columns.Item(10).ColumnWidth = $wideFileColWidth

Write-Output "Handback report generated"
